$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the balance typo on row 6 (1009 -> 1000)
$ws.Cells.Item(6, 6).Value = 1000

# Append 36 new "moses/bro" user rows (rows 16-51), cloned from row 6's
# text layout (Username/Password/ID/Email/Gender) so the cells keep their
# shared-string type + default style, with balance reset to 0 for each.
$ws.Range("A6:F6").Copy()
for ($r = 16; $r -le 51; $r++) {
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial()
    $ws.Cells.Item($r, 6).Value = 0
}
